# "version final sin errores"
#
# 1. Bump the ValueSet version number on the Metadata sheet (B3): 0.4.0 -> 0.7.0
# 2. Remove the "Jurisdiction" / "Chile" metadata row entirely (was row 11),
#    shifting all following rows up by one and shrinking the used range
#    from A1:B15 down to A1:B14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value.
$ws.Range("B3").Value = "0.7.0"

# Delete the entire "Jurisdiction" row (row 11), shifting subsequent rows up.
$ws.Rows.Item(11).Delete()
